# Insert a new weekly record at row 252, pushing the existing rows
# 252-371 down to 253-372 (new dimension becomes A1:R372).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 252, shifting rows down.
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with the new record's data.
$ws.Cells.Item(252, 1).Value  = 8
$ws.Cells.Item(252, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(252, 3).Value  = "Coquimbo"
$ws.Cells.Item(252, 4).Value  = 45205
$ws.Cells.Item(252, 5).Value  = 4
$ws.Cells.Item(252, 6).Value  = 100112037
$ws.Cells.Item(252, 7).Value  = "Cebollín"
$ws.Cells.Item(252, 8).Value  = "Sin especificar"
$ws.Cells.Item(252, 9).Value  = "Primera"
$ws.Cells.Item(252, 10).Value = 1400
$ws.Cells.Item(252, 11).Value = 1000
$ws.Cells.Item(252, 12).Value = 1200
$ws.Cells.Item(252, 13).Value = 1100
$ws.Cells.Item(252, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(252, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(252, 16).Value = 183
$ws.Cells.Item(252, 17).Value = 6
$ws.Cells.Item(252, 18).Value = "Hortaliza"
